$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A38").Value = "UserPageSizeDropdown"
$ws.Range("B38").Value = "//button[normalize-space()='10']"
$ws.Range("C38").Value = "By.xpath"
